$wb = $excel.ActiveWorkbook

# --- Add Denmark (copied from Turkey, placed right before Turkey) ---
$src = $wb.Worksheets.Item("Turkey")
$src.Copy($src, $null)
$denmark = $wb.Worksheets.Item("Turkey (2)")
$denmark.Name = "Denmark"
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").Value = "NGC-3446/T2004"
$denmark.Activate()
$denmark.Range("A1:XFD1048576").Select()

# --- Add Sweden (copied from Turkey, placed right before Turkey) ---
$src = $wb.Worksheets.Item("Turkey")
$src.Copy($src, $null)
$sweden = $wb.Worksheets.Item("Turkey (2)")
$sweden.Name = "Sweden"
$sweden.Range("B2").Value = "Sweden Market"
$sweden.Range("B4").Value = "NGC-3465/T2025"
$sweden.Activate()
$sweden.Range("A1:XFD1048576").Select()

# --- Add Norway (copied from Turkey, placed right before Turkey) ---
$src = $wb.Worksheets.Item("Turkey")
$src.Copy($src, $null)
$norway = $wb.Worksheets.Item("Turkey (2)")
$norway.Name = "Norway"
$norway.Range("B2").Value = "Norway Market"
$norway.Range("B4").Value = "NGC-3464/T1919"
$norway.Activate()
$norway.Range("B2:B4").Select()
